$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"

[void]$ws.Range("C1").Select()
